$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark: remove it from its old spot (the empty
#    paragraph right before "An additional multiple regression...") and
#    add it right before "Multiple regression was employed to relate"
#    (the first occurrence of that sentence).
# ---------------------------------------------------------------------
try {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
} catch {
}

$rngMR = $d.Content
$rngMR.Find.Execute("Multiple regression was employed to relate", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngMR.Collapse(1)
$d.Bookmarks.Add("_GoBack", $rngMR)

# ---------------------------------------------------------------------
# 2) Expand "... p < .001." into "... p < .001, 95% CI [-.01, .11]."
#    (only the FIRST occurrence of this sentence, which reports the R^2
#    for the emotional-exhaustion model) and push the sentence that used
#    to follow it directly ("Emotional exhaustion was positively related
#    to ...") into a new paragraph of its own.
# ---------------------------------------------------------------------
$statRange = $d.Content
$statRange.Find.Execute("F(2, 197) = 5.5, p < .001.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# statRange now covers exactly "F(2, 197) = 5.5, p < .001." - the very
# last character is the period we need to drop.
$periodRange = $d.Range($statRange.End - 1, $statRange.End)
$periodRange.Text = ""

$insPoint = $d.Range($periodRange.Start, $periodRange.Start)
$insPoint.InsertAfter(", 95% CI [-.01, .11].")

# The two spaces that used to separate the sentences become the new
# paragraph break.
$gap = $d.Range($insPoint.End, $insPoint.End + 2)
$gap.Text = [char]13

# ---------------------------------------------------------------------
# 3) Drop the stale "lastRenderedPageBreak" marker sitting on the
#    "Figure " run right before the Figure 1 caption. Word clears this
#    cached layout marker whenever that text is touched, so we delete
#    and retype it.
# ---------------------------------------------------------------------
$figRange = $d.Content
$figRange.Find.Execute("Figure 1.  Koleos", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$figWord = $d.Range($figRange.Start, $figRange.Start + 7)
$figWord.Text = ""
$figIns = $d.Range($figWord.Start, $figWord.Start)
$figIns.InsertAfter("Figure ")

Write-Output "done"
